# Applies the crypto price/volume refresh described by the commit
# "Updated cryptos list on Sat Jun 22 14:40:43 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.323.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.504.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.29%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.42%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.487"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.43%  "

# Row 9
$ws.Range("E9").Value = "  +0.88%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.389"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.77%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.103.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.12%  "

# Row 13
$ws.Range("E13").Value = "  +1.12%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.20%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.503.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.39%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.97%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.331.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.26%  "

# Row 18
$ws.Range("E18").Value = "  -3.25%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.38%  "

# Row 20
$ws.Range("E20").Value = "  -4.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "393.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.571"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.645.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.34%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "

# Row 26
$ws.Range("E26").Value = "  +2.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000114"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "

# Row 29
$ws.Range("E29").Value = "  +0.46%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "

# Row 31
$ws.Range("E31").Value = "  -1.11%  "

# Row 32
$ws.Range("E32").Value = "  -5.76%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.526.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.24%  "

# Row 34
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("E35").Value = "  +2.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.13%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.36%  "

# Row 38
$ws.Range("E38").Value = "  -1.03%  "

# Row 39
$ws.Range("E39").Value = "  -0.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.95%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0785"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.47%  "

# Row 42
$ws.Range("E42").Value = "  -1.34%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.18%  "

# Row 44
$ws.Range("E44").Value = "  +0.24%  "

# Row 45
$ws.Range("E45").Value = "  +0.66%  "

# Row 46
$ws.Range("E46").Value = "  -0.17%  "

# Row 47
$ws.Range("E47").Value = "  -4.27%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.30%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.468.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.22%  "

# Row 50
$ws.Range("E50").Value = "  -1.39%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.897"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
